# Refined metadata to be additional tab
#
# 1. Update the "panel_query_time" timestamps recorded in the `data` sheet
#    (column F) to reflect the re-run of the panel export.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name, id, version, version date, query time, request URL).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the recorded query timestamps on the data sheet ----------
$timeUpdates = @{
    2  = "2021-10-05 14:34:36.382991"
    3  = "2021-10-05 14:34:36.382999"
    4  = "2021-10-05 14:34:36.383002"
    5  = "2021-10-05 14:34:36.383005"
    6  = "2021-10-05 14:34:36.383008"
    7  = "2021-10-05 14:34:36.383011"
    8  = "2021-10-05 14:34:36.383013"
    9  = "2021-10-05 14:34:36.383015"
    10 = "2021-10-05 14:34:36.383018"
    11 = "2021-10-05 14:34:36.383021"
    12 = "2021-10-05 14:34:36.383023"
    13 = "2021-10-05 14:34:36.383026"
    14 = "2021-10-05 14:34:36.383028"
}

foreach ($row in $timeUpdates.Keys) {
    $dataSheet.Range("F$row").Value = $timeUpdates[$row]
}

# --- 2. Add the new "metadata" worksheet, placed after "data" ------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the header formatting (bold, centered, bordered) from the data sheet
# so the new sheet reuses the same cell style rather than minting a new one.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

$metaSheet.Application.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Mendelian susceptibility to Immune Disorders"
$metaSheet.Range("C2").Value = 228

# data_version is stored as literal text "0.13", not the number 0.13
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.13"

$metaSheet.Range("E2").Value = "2021-08-25T00:14:27.486981Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:36.379151"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/228/?format=json"

# Leave the user back on the first sheet, matching the original workbook view.
$dataSheet.Activate()

Write-Output "metadata tab added and data timestamps refreshed"
